# Weekly update: insert a new price record for "Camote" (Vega Modelo de
# Temuco) and push the existing history rows down by one row.
#
# The new record is inserted as row 77 (most recent date), and the rows
# that used to be 77-83 become 78-84 (their contents are unchanged by
# this insert - Excel shifts them down automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 77, shifting rows 77:83 -> 78:84.
$ws.Rows("77:77").Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "Vega Modelo de Temuco"
$ws.Range("C77").Value = "La Araucanía"
$ws.Range("D77").Value = 44706
$ws.Range("E77").Value = 9
$ws.Range("F77").Value = 100114002
$ws.Range("G77").Value = "Camote"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 50
$ws.Range("K77").Value = 18000
$ws.Range("L77").Value = 18000
$ws.Range("M77").Value = 18000
$ws.Range("N77").Value = "$/caja 15 kilos granel"
$ws.Range("O77").Value = "Perú"
$ws.Range("P77").Value = 1200
$ws.Range("Q77").Value = 15
$ws.Range("R77").Value = "Hortaliza"
